$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (style index) from H1 onto the new header cells
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

$ws.Range("I2").Value = 8
$ws.Range("J2").Value = 8

$ws.Range("I3").Value = 9
$ws.Range("J3").Value = 9
